# Auto-generated Excel COM-interop edit script
# Applies numeric cell-value updates, one deletion set, and one new-cell addition
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1110.1428
$ws.Range("I28").Value = 752
$ws.Range("J28").Value = 1826.4286
$ws.Range("K28").Value = 752
$ws.Range("L28").Value = 1826.4286
$ws.Range("M28").Value = -267
$ws.Range("N28").Value = -2796.4286
$ws.Range("H51").Value = 4779.9165
$ws.Range("I51").Value = 4279
$ws.Range("J51").Value = 4825.4546
$ws.Range("K51").Value = 4279
$ws.Range("L51").Value = 4825.4546
$ws.Range("M51").Value = -3795
$ws.Range("N51").Value = -5793.4546
$ws.Range("H70").Value = 1590.3077
$ws.Range("I70").Value = 1698.8
$ws.Range("J70").Value = 1522.5
$ws.Range("K70").Value = 5096.4
$ws.Range("L70").Value = 4567.5
$ws.Range("M70").Value = -4826.4
$ws.Range("N70").Value = -5107.5
$ws.Range("H73").Value = 1590.3077
$ws.Range("I73").Value = 1698.8
$ws.Range("J73").Value = 1522.5
$ws.Range("K73").Value = 5096.4
$ws.Range("L73").Value = 4567.5
$ws.Range("M73").Value = -4160.4
$ws.Range("N73").Value = -6439.5
$ws.Range("H76").Value = 10299.5625
$ws.Range("I76").Value = 7770.9287
$ws.Range("K76").Value = 7770.9287
$ws.Range("M76").Value = -7455.9287
$ws.Range("H79").Value = 10299.5625
$ws.Range("I79").Value = 7770.9287
$ws.Range("K79").Value = 7770.9287
$ws.Range("M79").Value = -6678.9287
$ws.Range("H97").Value = 2180.5833
$ws.Range("J97").Value = 2180.5833
$ws.Range("L97").Value = 6541.749899999999
$ws.Range("N97").Value = -7533.749899999999
$ws.Range("H98").Value = 3996.5
$ws.Range("I98").Value = 4132.636
$ws.Range("K98").Value = 4132.636
$ws.Range("M98").Value = -2634.636
$ws.Range("H106").Value = 58628.9
$ws.Range("I106").Value = 69761.25
$ws.Range("K106").Value = 69761.25
$ws.Range("M106").Value = -69130.25
$ws.Range("H112").Value = 1436.4546
$ws.Range("I112").Value = 359.5
$ws.Range("J112").Value = 1675.7778
$ws.Range("K112").Value = 1078.5
$ws.Range("L112").Value = 5027.3334
$ws.Range("M112").Value = 29.5
$ws.Range("N112").Value = -7243.3334
$ws.Range("H122").Value = 3996.5
$ws.Range("I122").Value = 4132.636
$ws.Range("K122").Value = 12397.908
$ws.Range("M122").Value = -9947.908000000001
$ws.Range("H128").Value = 89995
$ws.Range("J128").Value = 89995
$ws.Range("L128").Value = 89995
$ws.Range("N128").Value = -99955
$ws.Range("H138").Value = 1909.7391
$ws.Range("I138").Value = 971.1429000000001
$ws.Range("J138").Value = 2320.375
$ws.Range("K138").Value = 2913.4287
$ws.Range("L138").Value = 6961.125
$ws.Range("M138").Value = 2226.5713
$ws.Range("N138").Value = -17241.125
$ws.Range("H141").Value = 2706.4243
$ws.Range("I141").Value = 2781.1614
$ws.Range("K141").Value = 8343.484199999999
$ws.Range("M141").Value = -3163.484199999999

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("H13").Value = 25000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 25000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 25000
$ws.Range("N13").Value = -25288
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("H32").Value = 4071.6592
$ws.Range("I32").Value = 2064.353
$ws.Range("K32").Value = 2064.353
$ws.Range("M32").Value = -1777.353
$ws.Range("H45").Value = 1630.625
$ws.Range("I45").Value = 1439.7142
$ws.Range("J45").Value = 1995.091
$ws.Range("K45").Value = 1439.7142
$ws.Range("L45").Value = 1995.091
$ws.Range("M45").Value = -1062.7142
$ws.Range("N45").Value = -2749.091
$ws.Range("H61").Value = 4252.8076
$ws.Range("I61").Value = 3934
$ws.Range("J61").Value = 4486.6
$ws.Range("K61").Value = 3934
$ws.Range("L61").Value = 4486.6
$ws.Range("M61").Value = -3722
$ws.Range("N61").Value = -4910.6
$ws.Range("H88").Value = 2076.1853
$ws.Range("I88").Value = 1848.8889
$ws.Range("J88").Value = 2189.8333
$ws.Range("K88").Value = 1848.8889
$ws.Range("L88").Value = 2189.8333
$ws.Range("M88").Value = -1442.8889
$ws.Range("N88").Value = -3001.8333
$ws.Range("H91").Value = 2076.1853
$ws.Range("I91").Value = 1848.8889
$ws.Range("J91").Value = 2189.8333
$ws.Range("K91").Value = 1848.8889
$ws.Range("L91").Value = 2189.8333
$ws.Range("M91").Value = -444.8888999999999
$ws.Range("N91").Value = -4997.8333
$ws.Range("H136").Value = 4252.8076
$ws.Range("I136").Value = 3934
$ws.Range("J136").Value = 4486.6
$ws.Range("K136").Value = 11802
$ws.Range("L136").Value = 13459.8
$ws.Range("M136").Value = -9252
$ws.Range("N136").Value = -18559.8
$ws.Range("M11").ClearContents()
$ws.Range("M13").ClearContents()
$ws.Range("N23").ClearContents()

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2512
$ws.Range("I20").Value = 3293
$ws.Range("K20").Value = 3293
$ws.Range("M20").Value = -3046
$ws.Range("H86").Value = 1720.3125
$ws.Range("J86").Value = 2006
$ws.Range("L86").Value = 2006
$ws.Range("N86").Value = -4252
$ws.Range("H89").Value = 1720.3125
$ws.Range("J89").Value = 2006
$ws.Range("L89").Value = 10030
$ws.Range("N89").Value = -21262
$ws.Range("H107").Value = 2686.9312
$ws.Range("I107").Value = 1813.1428
$ws.Range("K107").Value = 1813.1428
$ws.Range("M107").Value = 106.8571999999999

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 11263.223
$ws.Range("I105").Value = 13202.6
$ws.Range("J105").Value = 1566.3334
$ws.Range("K105").Value = 13202.6
$ws.Range("L105").Value = 1566.3334
$ws.Range("M105").Value = -11455.6
$ws.Range("N105").Value = -5060.3334
$ws.Range("H134").Value = 2015.6
$ws.Range("I134").Value = 1884.7142
$ws.Range("J134").Value = 2702.75
$ws.Range("K134").Value = 5654.142599999999
$ws.Range("L134").Value = 8108.25
$ws.Range("M134").Value = -3119.142599999999
$ws.Range("N134").Value = -13178.25

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 323.33334
$ws.Range("J92").Value = 335.5
$ws.Range("L92").Value = 1006.5
$ws.Range("N92").Value = -3502.5
$ws.Range("H93").Value = 700
$ws.Range("I93").Value = 700
$ws.Range("K93").Value = 2100
$ws.Range("M93").Value = -228
$ws.Range("H99").Value = 2475
$ws.Range("I99").Value = 2475
$ws.Range("K99").Value = 7425
$ws.Range("M99").Value = -5179
$ws.Range("H100").Value = 78942.5
$ws.Range("J100").Value = 78942.5
$ws.Range("L100").Value = 236827.5
$ws.Range("N100").Value = -238449.5
$ws.Range("H137").Value = 2589.1538
$ws.Range("J137").Value = 2682.9
$ws.Range("L137").Value = 8048.700000000001
$ws.Range("N137").Value = -18248.7

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 69904.39999999999
$ws.Range("J135").Value = 74880.5
$ws.Range("L135").Value = 74880.5
$ws.Range("N135").Value = -85020.5

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5978.1577
$ws.Range("I7").Value = 3069.6
$ws.Range("J7").Value = 9209.888999999999
$ws.Range("K7").Value = 3069.6
$ws.Range("L7").Value = 9209.888999999999
$ws.Range("M7").Value = -2957.6
$ws.Range("N7").Value = -9433.888999999999
$ws.Range("H22").Value = 3897.4167
$ws.Range("I22").Value = 1734.4
$ws.Range("K22").Value = 1734.4
$ws.Range("M22").Value = -1439.4
$ws.Range("H27").Value = 3897.4167
$ws.Range("I27").Value = 1734.4
$ws.Range("K27").Value = 1734.4
$ws.Range("M27").Value = -1627.4
$ws.Range("H68").Value = 1414.7742
$ws.Range("I68").Value = 1868.7693
$ws.Range("K68").Value = 1868.7693
$ws.Range("M68").Value = -1119.7693
$ws.Range("H71").Value = 1414.7742
$ws.Range("I71").Value = 1868.7693
$ws.Range("K71").Value = 9343.8465
$ws.Range("M71").Value = -5599.8465
$ws.Range("H126").Value = 5978.1577
$ws.Range("I126").Value = 3069.6
$ws.Range("J126").Value = 9209.888999999999
$ws.Range("K126").Value = 9208.799999999999
$ws.Range("L126").Value = 27629.667
$ws.Range("M126").Value = -6738.799999999999
$ws.Range("N126").Value = -32569.667
$ws.Range("H136").Value = 4409.1665
$ws.Range("I136").Value = 3909.5454
$ws.Range("K136").Value = 11728.6362
$ws.Range("M136").Value = -9178.636200000001

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 19850
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("H126").Value = 4454.1816
$ws.Range("I126").Value = 4454.1816
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 13362.5448
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -10892.5448
$ws.Range("N3").ClearContents()
$ws.Range("N126").ClearContents()
